$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-06 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-07 Monday", 2) | Out-Null
$d.Content.Find.Execute("246÷9=27, 3", $true, $false, $false, $false, $false, $true, 1, $false, "294÷5=58, 4", 2) | Out-Null
$d.Content.Find.Execute("165÷2=82, 1", $true, $false, $false, $false, $false, $true, 1, $false, "439÷8=54, 7", 2) | Out-Null
$d.Content.Find.Execute("523÷3=174, 1", $true, $false, $false, $false, $false, $true, 1, $false, "818÷6=136, 2", 2) | Out-Null
$d.Content.Find.Execute("321÷4=80, 1", $true, $false, $false, $false, $false, $true, 1, $false, "299÷6=49, 5", 2) | Out-Null
$d.Content.Find.Execute("134÷2=67, 0", $true, $false, $false, $false, $false, $true, 1, $false, "997÷3=332, 1", 2) | Out-Null
$d.Content.Find.Execute("582÷2=291, 0", $true, $false, $false, $false, $false, $true, 1, $false, "940÷4=235, 0", 2) | Out-Null
$d.Content.Find.Execute("885÷6=147, 3", $true, $false, $false, $false, $false, $true, 1, $false, "672÷2=336, 0", 2) | Out-Null
$d.Content.Find.Execute("596÷9=66, 2", $true, $false, $false, $false, $false, $true, 1, $false, "104÷3=34, 2", 2) | Out-Null
$d.Content.Find.Execute("331÷7=47, 2", $true, $false, $false, $false, $false, $true, 1, $false, "874÷4=218, 2", 2) | Out-Null
$d.Content.Find.Execute("783÷7=111, 6", $true, $false, $false, $false, $false, $true, 1, $false, "776÷6=129, 2", 2) | Out-Null
$d.Content.Find.Execute("214÷7=30, 4", $true, $false, $false, $false, $false, $true, 1, $false, "919÷4=229, 3", 2) | Out-Null
$d.Content.Find.Execute("841÷6=140, 1", $true, $false, $false, $false, $false, $true, 1, $false, "838÷4=209, 2", 2) | Out-Null
$d.Content.Find.Execute("690÷9=76, 6", $true, $false, $false, $false, $false, $true, 1, $false, "974÷6=162, 2", 2) | Out-Null
$d.Content.Find.Execute("227÷8=28, 3", $true, $false, $false, $false, $false, $true, 1, $false, "702÷6=117, 0", 2) | Out-Null
$d.Content.Find.Execute("733÷2=366, 1", $true, $false, $false, $false, $false, $true, 1, $false, "771÷4=192, 3", 2) | Out-Null
$d.Content.Find.Execute("696÷3=232, 0", $true, $false, $false, $false, $false, $true, 1, $false, "324÷3=108, 0", 2) | Out-Null
$d.Content.Find.Execute("505÷9=56, 1", $true, $false, $false, $false, $false, $true, 1, $false, "346÷5=69, 1", 2) | Out-Null
$d.Content.Find.Execute("263÷2=131, 1", $true, $false, $false, $false, $false, $true, 1, $false, "786÷5=157, 1", 2) | Out-Null
$d.Content.Find.Execute("829÷6=138, 1", $true, $false, $false, $false, $false, $true, 1, $false, "108÷2=54, 0", 2) | Out-Null
$d.Content.Find.Execute("790÷2=395, 0", $true, $false, $false, $false, $false, $true, 1, $false, "209÷4=52, 1", 2) | Out-Null
$d.Content.Find.Execute("559÷5=111, 4", $true, $false, $false, $false, $false, $true, 1, $false, "306÷2=153, 0", 2) | Out-Null
$d.Content.Find.Execute("253÷8=31, 5", $true, $false, $false, $false, $false, $true, 1, $false, "229÷4=57, 1", 2) | Out-Null
$d.Content.Find.Execute("576÷6=96, 0", $true, $false, $false, $false, $false, $true, 1, $false, "634÷7=90, 4", 2) | Out-Null
$d.Content.Find.Execute("408÷4=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "780÷9=86, 6", 2) | Out-Null
$d.Content.Find.Execute("245÷8=30, 5", $true, $false, $false, $false, $false, $true, 1, $false, "516÷4=129, 0", 2) | Out-Null
